# grump_GA10_cmap_submission.xlsx -- "Grump metadata. Longhurst provinces and
# Seasons added."
#
# The vars_meta_data sheet documents one variable per row (columns:
# var_short_name, var_long_name, var_sensor, var_unit, var_spatial_res,
# var_temporal_res, var_discipline, visualize). Rows 41-43 were blank
# placeholder rows; this adds three new variable-documentation rows describing
# the Longhurst province (long + short code) and Season fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "vars_meta_data" is the active/selected sheet in this workbook

# Row 40 ("DCM") is the last populated data row and already carries the
# per-column formatting (style) used throughout the table. Duplicating it via
# Copy+Insert (rather than just typing into the blank rows) brings that
# formatting along for the new rows, matching how the rest of the sheet looks.
$ws.Rows.Item(40).Copy()
$ws.Rows.Item(41).Insert()
$ws.Rows.Item(40).Copy()
$ws.Rows.Item(42).Insert()
$ws.Rows.Item(40).Copy()
$ws.Rows.Item(43).Insert()

# The Insert() calls above pushed the (already blank) former rows 41-43 down
# to 44-46 -- remove that now-duplicated blank padding so row numbering/count
# ends up unchanged.
$ws.Rows.Item(44).Delete()
$ws.Rows.Item(44).Delete()
$ws.Rows.Item(44).Delete()

# Restore the row height on the three new rows (13.5pt, matching the rest of
# the sheet) since Insert() doesn't always preserve it.
$ws.Rows.Item(41).RowHeight = 13.5
$ws.Rows.Item(42).RowHeight = 13.5
$ws.Rows.Item(43).RowHeight = 13.5

# Row 41: Longhurst province (long form)
$ws.Range("A41").Value = "Longhurst_Long"
$ws.Range("B41").Value = "Longhurst province sample was taken in."
$ws.Range("C41").Value = "NA"
$ws.Range("D41").Value = "NA"
$ws.Range("E41").Value = "Irregular"
$ws.Range("F41").Value = "Irregular"
$ws.Range("G41").Value = "Biology"
$ws.Range("H41").Value = 1

# Row 42: Longhurst province (shortened code)
$ws.Range("A42").Value = "Longhurst_Short"
$ws.Range("B42").Value = "Longhurst province sample was taken in, shortened code."
$ws.Range("C42").Value = "NA"
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "Irregular"
$ws.Range("F42").Value = "Irregular"
$ws.Range("G42").Value = "Biology"
$ws.Range("H42").Value = 1

# Row 43: Season
$ws.Range("A43").Value = "Season"
$ws.Range("B43").Value = "Season sample was taken in."
$ws.Range("C43").Value = "NA"
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = "Irregular"
$ws.Range("F43").Value = "Irregular"
$ws.Range("G43").Value = "Biology"
$ws.Range("H43").Value = 1

# Leave the selection on the newly-added block, matching the author's final
# on-screen state after typing the new rows in.
$ws.Activate()
$ws.Range("A41:H43").Select()
